$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: new assignment entry (Test Student, due 01/05/1235)
$ws.Range("A2").Formula = "'01/05/1235"
$ws.Range("A2").ClearFormats()
$ws.Range("B2").Value = "Test"
$ws.Range("C2").Value = "Student"
$ws.Range("D2").Formula = "'1234"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Formula = "'0"
$ws.Range("E2").ClearFormats()

# Row 3: new assignment entry (Fisrt Last, due 01/10/1001)
$ws.Range("A3").Formula = "'01/10/1001"
$ws.Range("A3").ClearFormats()
$ws.Range("B3").Value = "Fisrt"
$ws.Range("C3").Value = "Last"
$ws.Range("D3").Formula = "'1234"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Formula = "'0"
$ws.Range("E3").ClearFormats()
